$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-28 all hold the date serial 45548 (2024-09-13);
# bump each one by a day to 45549 (2024-09-14), matching the source diff.
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45548) {
        $cell.Value2 = 45549
    }
}
